# Apply updated measurement values to week 1 - week 4 columns (F:I) for rows 2-37.
# Values reflect the canonical OOXML diff for this workbook revision
# (sex/status categories expanded and analysis values recomputed; "talk" note added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row's new values for columns F, G, H, I (in order).
$newValues = @(
    @(1.62, 6.89, 7.8, 12.82),
    @(5.25, 10.18, 9.92, 12.49),
    @(5.33, 9.369999999999999, 11.21, 15.51),
    @(7.85, 10.41, 12.72, 14.27),
    @(8.26, 9.43, 11.89, 14.92),
    @(11.06, 9.25, 13.36, 18.77),
    @(4.33, 10.49, 8.880000000000001, 12.1),
    @(2.43, 8.44, 12.88, 1205),
    @(8.48, 8.710000000000001, 7.73, 16.56),
    @(9.15, 10.84, 14.68, 1534),
    @(7.06, 12.4, 14.42, 17.89),
    @(9.390000000000001, 13.28, 18.08, 17.92),
    @(3, 12.08, 9.859999999999999, 8.9),
    @(6.78, 11.79, 9.52, 11.65),
    @(5.48, 10.91, 13.01, 14.55),
    @(6.17, 10.86, 15.34, 16.11),
    @(732, 12.31, 11.96, 17.6),
    @(9.57, 13.03, 13.89, 15.98),
    @(9.220000000000001, 6.19, 7.5, 13.81),
    @(3.99, 5.39, 13.5, 11.74),
    @(8.58, 8.31, 13.13, 14.09),
    @(10.35, 8.119999999999999, 15.93, 14.6),
    @(7.82, 11.89, 15.29, 16.56),
    @(5.56, 13.57, 16.08, 15.79),
    @(7.24, 9.65, 11.78, 16.04),
    @(6.6, 3.8, 10.94, 12.75),
    @(5.86, 10.01, 12.04, 15.01),
    @(12.63, 9.119999999999999, 15.52, 17.22),
    @(8.289999999999999, 14.22, 16.68, 17.97),
    @(12.26, 10.65, 14.82, 15.41),
    @(5.36, 5.24, 9.07, 12.64),
    @(5, 3.88, 10.53, 14.36),
    @(3.59, 10.49, 9.02, 14.22),
    @(6.48, 9.92, 1288, 18.83),
    @(7.05, 1116, 9.91, 17.14),
    @(9.550000000000001, 12.78, 12.12, 18.48)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 6).Value = $rowValues[0]
    $ws.Cells.Item($r, 7).Value = $rowValues[1]
    $ws.Cells.Item($r, 8).Value = $rowValues[2]
    $ws.Cells.Item($r, 9).Value = $rowValues[3]
}
